# Rerun: rename the per-sheet "summ########" identifiers to a fresh set,
# keeping sheet order / positions unchanged.
$wb = $excel.ActiveWorkbook

$newNames = @(
    "summ27682872",
    "summ52280198",
    "summ17759501",
    "summ48329543",
    "summ19134248",
    "summ45819104",
    "summ13670133",
    "summ40812309",
    "summ08065431"
)

for ($i = 1; $i -le $wb.Worksheets.Count; $i++) {
    $wb.Worksheets.Item($i).Name = $newNames[$i - 1]
}
